$d = $word.ActiveDocument

# Cursor-based sequential replace: always search forward from $cursor
# so that ambiguous/duplicate strings only match the next occurrence in
# document order (matching the order the diff hunks appear in the file).
$docEnd = $d.Content.End

function Replace-Next([string]$find, [string]$repl) {
    $r = $d.Range($cursor, $docEnd)
    $found = $r.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 1)
    if (-not $found) {
        throw "Find failed for: $find"
    }
    $script:cursor = $r.End
    Write-Output ("Replaced [" + $find.Substring(0, [Math]::Min(30, $find.Length)) + "...] OK, cursor=" + $script:cursor)
}

# The document starts with a navigation line "English / Portuguese / French /
# Thai / Vietnamese / Spanish" whose items are hyperlinks to each language's
# section below (e.g. the first "English" is a hyperlink, not the section
# heading). Skip past that line first so the first real replacement below
# lands on the "English" section heading, not the nav-bar hyperlink.
$cursor = 0
$navRange = $d.Range(0, $docEnd)
$navFound = $navRange.Find.Execute("English / Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $navFound) {
    throw "Could not locate the language navigation line"
}
$cursor = $navRange.End

# --- English section (translated to Thai) ---
Replace-Next 'English' 'ภาษาอังกฤษ'
Replace-Next 'Brief' 'บทย่อ'
Replace-Next 'An email to partners in the the target country to invite them for a one-day seminar. It will be sent via customer.io' 'อีเมล์ถึงพันธมิตรในประเทศเป้าหมายเพื่อเชิญพวกเขาเข้าร่วมสัมมนาที่จัดขึ้นภายในวันเดียว โดยมันจะถูกส่งผ่านทาง customer.io'
Replace-Next 'Target audience' 'กลุ่มเป้าหมาย'
Replace-Next 'Partners in the target country' 'พันธมิตรหุ้นส่วนในประเทศเป้าหมาย'
Replace-Next 'Subject line' 'หัวเรื่อง'
Replace-Next ': Meet our team in [CITY] | [DATE] ' ': พบกับทีมของเราได้ใน [CITY] | [DATE] '
Replace-Next 'You’re invited to our Deriv Partner Seminar' 'คุณได้รับเชิญให้เข้าร่วมงานสัมมนาหุ้นส่วน Deriv'
Replace-Next 'Dear [PARTNER NAME], ' 'เรียนคุณ [PARTNER NAME] '
Replace-Next 'We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!' 'เรารู้สึกตื่นเต้นที่จะแจ้งให้คุณทราบว่า ทีมพันธมิตร Deriv จะเยือน [CITY] ในเดือน [MONTH] เพื่อพบปะกับคุณผู้ซึ่งเป็นพันธมิตรที่มีค่าของเรา!'
Replace-Next 'Your country manager will inform you about the exact location by [DATE]' 'ผู้จัดการประเทศของคุณจะแจ้งให้คุณทราบเกี่ยวกับสถานที่จัดงานภายในวันที่ [DATE]'
Replace-Next 'In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. ' 'ในการสัมมนาหนึ่งวันนี้ เราจะให้การสนับสนุนด้านเทคนิคและการตลาด เปิดโอกาสให้คุณได้สร้างเครือข่ายกับพันธมิตรรายอื่นๆ ในระหว่างการรับประทานอาหารกลางวันแสนอร่อย รวมทั้งรับฟังความคิดเห็นของคุณเกี่ยวกับโครงการหุ้นส่วนพันธมิตรต่างๆ ของเรา นี่เป็นโอกาสของคุณที่จะแสดงความเห็นซึ่งจะช่วยให้เราวางแผนความพยายามสนับสนุนคุณในอนาคตให้ดียิ่งขึ้น '
Replace-Next 'Please RSVP by submitting the registration form by ' 'โปรดตอบกลับ RSVP โดยกรอกแบบฟอร์มลงทะเบียนมาให้เราภายในวันที่ '
Replace-Next '. Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!' ' โปรดทราบว่า การเข้าร่วมประชุมจะพิจารณายืนยันไปตามลำดับใครมาก่อนได้ก่อน เราหวังว่าจะได้พบเจอคุณที่นั่น!'
Replace-Next 'Send my details' 'ส่งรายละเอียดของฉัน'
Replace-Next 'If you have any questions, please contact us via ' 'หากคุณมีคำถามใดๆ กรุณาติดต่อเราผ่านทาง '
Replace-Next 'live chat' 'แชทสด'
Replace-Next ' or ' ' หรือทาง '
Replace-Next '. / If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). ' ' / หากคุณมีคำถามใดๆ โปรดติดต่อผู้จัดการประจำประเทศของคุณซึ่งได้แก่ [NAME] ที่ [EMAIL ADDRESS] หรือ [WHATSAPP NO] (WhatsApp) '

# --- Second 'Send my details' button (French section, still English in source) ---
Replace-Next 'Send my details' 'ส่งรายละเอียดของฉัน'

# --- Minor corrections within the existing Thai section ---
Replace-Next ': พบกับทีมงานของเราที่ [CITY] | [DATE]' ': พบกับทีมงานของเราได้ที่ [CITY] | [DATE]'
Replace-Next 'Deriv Partner Seminar ของเรา' 'สัมมนาพันธมิตรหุ้นส่วน Deriv ของเรา'

Write-Output "DONE"
